$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Board_Layout")

# Fill column Z (rows 1-25) with values 0..24
for ($row = 1; $row -le 25; $row++) {
    $ws.Cells.Item($row, 26).Value = $row - 1
}

# Fill row 26 (columns A-Y, i.e. 1-25) with values 0..24
for ($col = 1; $col -le 25; $col++) {
    $ws.Cells.Item(26, $col).Value = $col - 1
}

# Update the active cell selection to H6
$ws.Range("H6").Select()
